$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.8125;              "C2" = 0.7959183673469388;  "D2" = 0.8041237113402061;  "E2" = 49
    "B3" = 0.7777777777777778;  "C3" = 0.7954545454545454;  "D3" = 0.7865168539325843;  "E3" = 44
    "B4" = 0.7956989247311828;  "C4" = 0.7956989247311828;  "D4" = 0.7956989247311828;  "E4" = 0.7956989247311828
    "B5" = 0.7951388888888888;  "C5" = 0.7956864564007421;  "D5" = 0.7953202826363952
    "B6" = 0.7960722819593788;  "C6" = 0.7956989247311828;  "D6" = 0.7957935852548795
    "B7" = 0.8627450980392157;  "C7" = 0.8979591836734694;  "D7" = 0.8799999999999999;  "E7" = 49
    "B8" = 0.8809523809523809;  "C8" = 0.8409090909090909;  "D8" = 0.8604651162790699;  "E8" = 44
    "B9" = 0.8709677419354839;  "C9" = 0.8709677419354839;  "D9" = 0.8709677419354839;  "E9" = 0.8709677419354839
    "B10" = 0.8718487394957983; "C10" = 0.8694341372912802; "D10" = 0.8702325581395349
    "B11" = 0.8713592964067348; "C11" = 0.8709677419354839; "D11" = 0.8707576894223557
    "B12" = 0.7450980392156863; "C12" = 0.7755102040816326; "D12" = 0.76;                "E12" = 49
    "B13" = 0.7380952380952381; "C13" = 0.7045454545454546; "D13" = 0.7209302325581395;  "E13" = 44
    "B14" = 0.7419354838709677; "C14" = 0.7419354838709677; "D14" = 0.7419354838709677;  "E14" = 0.7419354838709677
    "B15" = 0.7415966386554622; "C15" = 0.7400278293135436; "D15" = 0.7404651162790697
    "B16" = 0.7417848859974097; "C16" = 0.7419354838709677; "D16" = 0.7415153788447112
    "B17" = 0.8181818181818182; "C17" = 0.9183673469387755; "D17" = 0.8653846153846154;  "E17" = 49
    "B18" = 0.8947368421052632; "C18" = 0.7727272727272727; "D18" = 0.8292682926829269;  "E18" = 44
    "B19" = 0.8494623655913979; "C19" = 0.8494623655913979; "D19" = 0.8494623655913979;  "E19" = 0.8494623655913979
    "B20" = 0.8564593301435407; "C20" = 0.8455473098330242; "D20" = 0.8473264540337712
    "B21" = 0.8544013993929106; "C21" = 0.8494623655913979; "D21" = 0.8482973229236015
    "B22" = 1; "C22" = 1; "D22" = 1; "E22" = 49
    "B23" = 1; "C23" = 1; "D23" = 1; "E23" = 44
    "B24" = 1; "C24" = 1; "D24" = 1; "E24" = 1
    "B25" = 1; "C25" = 1; "D25" = 1
    "B26" = 1; "C26" = 1; "D26" = 1
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
